# Update pl_mw.xlsx results for Case_4_153 (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.389736918993435
$ws.Range("C2").Value = 0.08370356941257739
$ws.Range("D2").Value = 0.07561659067076221
$ws.Range("E2").Value = 0.03554583425865143
$ws.Range("G2").Value = 0.002673195386167774
$ws.Range("K2").Value = 1.953148984501354
$ws.Range("L2").Value = 0.2007485387920624
$ws.Range("M2").Value = 0.4458007969451216
$ws.Range("N2").Value = 5.950383911107537
$ws.Range("B3").Value = 2.33813304586721
$ws.Range("C3").Value = 0.07503612207599986
$ws.Range("D3").Value = 0.06904561120113328
$ws.Range("E3").Value = 0.03563475596782428
$ws.Range("G3").Value = 0.002678928241498973
$ws.Range("K3").Value = 1.89533729178541
$ws.Range("L3").Value = 0.1986554679801387
$ws.Range("M3").Value = 0.4375680211641395
$ws.Range("N3").Value = 5.870573711976476
$ws.Range("B4").Value = 2.307828134379946
$ws.Range("C4").Value = 0.0697548868561455
$ws.Range("D4").Value = 0.06505379672303491
$ws.Range("E4").Value = 0.03569518316566667
$ws.Range("G4").Value = 0.002682631639817412
$ws.Range("K4").Value = 1.861035546782773
$ws.Range("L4").Value = 0.1974714040244905
$ws.Range("M4").Value = 0.4327632486765154
$ws.Range("N4").Value = 5.821985800861881
$ws.Range("B5").Value = 2.295824721582647
$ws.Range("C5").Value = 0.06761270454708779
$ws.Range("D5").Value = 0.06343770596136267
$ws.Range("E5").Value = 0.03572127683829729
$ws.Range("G5").Value = 0.002684187085083702
$ws.Range("K5").Value = 1.847356551628366
$ws.Range("L5").Value = 0.197014284341634
$ws.Range("M5").Value = 0.4308680349132743
$ws.Range("N5").Value = 5.802289333108916
$ws.Range("B6").Value = 2.293852439587397
$ws.Range("C6").Value = 0.06725759139654031
$ws.Range("D6").Value = 0.06316999141547797
$ws.Range("E6").Value = 0.03572569850287488
$ws.Range("G6").Value = 0.002684448165328138
$ws.Range("K6").Value = 1.845103203570119
$ws.Range("L6").Value = 0.1969399131561786
$ws.Range("M6").Value = 0.4305571232401633
$ws.Range("N6").Value = 5.799024969604574
$ws.Range("B7").Value = 2.30766485211916
$ws.Range("C7").Value = 0.06972595661684977
$ws.Range("D7").Value = 0.06503195881538204
$ws.Range("E7").Value = 0.03569552912095109
$ws.Range("G7").Value = 0.002682652429522037
$ws.Range("K7").Value = 1.860849857036357
$ws.Range("L7").Value = 0.1974651363509281
$ws.Range("M7").Value = 0.432737435192287
$ws.Range("N7").Value = 5.821719749820147
$ws.Range("B8").Value = 2.371656870024481
$ws.Range("C8").Value = 0.08070645163024892
$ws.Range("D8").Value = 0.07334193951702161
$ws.Range("E8").Value = 0.03557528667474008
$ws.Range("G8").Value = 0.002675134107496053
$ws.Range("K8").Value = 1.932966757042891
$ws.Range("L8").Value = 0.2000058508010767
$ws.Range("M8").Value = 0.4429101201266832
$ws.Range("N8").Value = 5.922778366148975
$ws.Range("B9").Value = 2.508148577949214
$ws.Range("C9").Value = 0.1025737716745141
$ws.Range("D9").Value = 0.08998522577051915
$ws.Range("E9").Value = 0.03538559859513191
$ws.Range("G9").Value = 0.002661838495837369
$ws.Range("K9").Value = 2.083934028915451
$ws.Range("L9").Value = 0.2057919383245732
$ws.Range("M9").Value = 0.4648519769271431
$ws.Range("N9").Value = 6.124312908559091
$ws.Range("B10").Value = 2.615226189492887
$ws.Range("C10").Value = 0.118862458152762
$ws.Range("D10").Value = 0.1024371027844211
$ws.Range("E10").Value = 0.03527415888064933
$ws.Range("G10").Value = 0.002652942421106942
$ws.Range("K10").Value = 2.200775532140597
$ws.Range("L10").Value = 0.210536081609817
$ws.Range("M10").Value = 0.4822016605456341
$ws.Range("N10").Value = 6.274525286758774
$ws.Range("B11").Value = 2.665435113588558
$ws.Range("C11").Value = 0.1263251699834598
$ws.Range("D11").Value = 0.1081531838116376
$ws.Range("E11").Value = 0.03522948820165706
$ws.Range("G11").Value = 0.002649082534764346
$ws.Range("K11").Value = 2.255240794324379
$ws.Range("L11").Value = 0.2128021377048981
$ws.Range("M11").Value = 0.4903645726283798
$ws.Range("N11").Value = 6.343348548187691
$ws.Range("B12").Value = 2.684664879700676
$ws.Range("C12").Value = 0.1291590240153369
$ws.Range("D12").Value = 0.1103253425359298
$ws.Range("E12").Value = 0.0352134358023739
$ws.Range("G12").Value = 0.002647647613893989
$ws.Range("K12").Value = 2.276055969001902
$ws.Range("L12").Value = 0.2136758015808056
$ws.Range("M12").Value = 0.4934947673597421
$ws.Range("N12").Value = 6.369482178253065
$ws.Range("B13").Value = 2.680513746418512
$ws.Range("C13").Value = 0.1285483475803915
$ws.Range("D13").Value = 0.1098571883467514
$ws.Range("E13").Value = 0.03521685461687518
$ws.Range("G13").Value = 0.002647955463295787
$ws.Range("K13").Value = 2.271564558277703
$ws.Range("L13").Value = 0.2134869498702727
$ws.Range("M13").Value = 0.4928188829296687
$ws.Range("N13").Value = 6.36385062342265
$ws.Range("B14").Value = 2.667012807233277
$ws.Range("C14").Value = 0.1265581531795021
$ws.Range("D14").Value = 0.1083317352358648
$ws.Range("E14").Value = 0.03522815027271653
$ws.Range("G14").Value = 0.002648963948197473
$ws.Range("K14").Value = 2.256949447502222
$ws.Range("L14").Value = 0.2128737025325194
$ws.Range("M14").Value = 0.4906213113835278
$ws.Range("N14").Value = 6.345497132879927
$ws.Range("B15").Value = 2.658771354191799
$ws.Range("C15").Value = 0.1253401371195935
$ws.Range("D15").Value = 0.1073983464985702
$ws.Range("E15").Value = 0.0352351815486287
$ws.Range("G15").Value = 0.002649585151116675
$ws.Range("K15").Value = 2.248022105578571
$ws.Range("L15").Value = 0.2125000980683609
$ws.Range("M15").Value = 0.4892803299685653
$ws.Range("N15").Value = 6.334264469315315
$ws.Range("B16").Value = 2.611975180484535
$ws.Range("C16").Value = 0.1183758407680955
$ws.Range("D16").Value = 0.1020645978117045
$ws.Range("E16").Value = 0.03527719918807382
$ws.Range("G16").Value = 0.002653198423314312
$ws.Range("K16").Value = 2.197242655878256
$ws.Range("L16").Value = 0.2103901644333348
$ws.Range("M16").Value = 0.4816736530517289
$ws.Range("N16").Value = 6.270037502089423
$ws.Range("B17").Value = 2.583651970175424
$ws.Range("C17").Value = 0.1141172185017183
$ws.Range("D17").Value = 0.09880587023977228
$ws.Range("E17").Value = 0.03530451647660193
$ws.Range("G17").Value = 0.002655462829824677
$ws.Range("K17").Value = 2.166428567197158
$ws.Range("L17").Value = 0.2091234567013203
$ws.Range("M17").Value = 0.4770766035506853
$ws.Range("N17").Value = 6.230762892733907
$ws.Range("B18").Value = 2.56750220655988
$ws.Range("C18").Value = 0.1116727344983417
$ws.Range("D18").Value = 0.09693638820826322
$ws.Range("E18").Value = 0.0353207958104389
$ws.Range("G18").Value = 0.002656782865490187
$ws.Range("K18").Value = 2.148828658080618
$ws.Range("L18").Value = 0.2084050355837519
$ws.Range("M18").Value = 0.474457947204435
$ws.Range("N18").Value = 6.208219311395226
$ws.Range("B19").Value = 2.562058354172962
$ws.Range("C19").Value = 0.1108459181596686
$ws.Range("D19").Value = 0.09630424176229724
$ws.Range("E19").Value = 0.03532640520865549
$ws.Range("G19").Value = 0.002657232836101474
$ws.Range("K19").Value = 2.142890804610204
$ws.Range("L19").Value = 0.208163533634405
$ws.Range("M19").Value = 0.4735756797717627
$ws.Range("N19").Value = 6.200594330586796
$ws.Range("B20").Value = 2.586652422846441
$ws.Range("C20").Value = 0.114570040463235
$ws.Range("D20").Value = 0.09915226394792853
$ws.Range("E20").Value = 0.03530154982696887
$ws.Range("G20").Value = 0.002655219958651605
$ws.Range("K20").Value = 2.169695983003123
$ws.Range("L20").Value = 0.2092572485549482
$ws.Range("M20").Value = 0.4775633325334638
$ws.Range("N20").Value = 6.234938958741139
$ws.Range("B21").Value = 2.670972468342939
$ws.Range("C21").Value = 0.1271425052208599
$ws.Range("D21").Value = 0.1087795900499913
$ws.Range("E21").Value = 0.03522480905390291
$ws.Range("G21").Value = 0.002648667007269635
$ws.Range("K21").Value = 2.261237081431716
$ws.Range("L21").Value = 0.2130534055950193
$ws.Range("M21").Value = 0.49126572936418
$ws.Range("N21").Value = 6.350886041199658
$ws.Range("B22").Value = 2.727344384387834
$ws.Range("C22").Value = 0.135405441825867
$ws.Range("D22").Value = 0.1151160024510602
$ws.Range("E22").Value = 0.03517968593118947
$ws.Range("G22").Value = 0.002644540032090725
$ws.Range("K22").Value = 2.322174770470895
$ws.Range("L22").Value = 0.2156251145524521
$ws.Range("M22").Value = 0.500448899390733
$ws.Range("N22").Value = 6.427082952955857
$ws.Range("B23").Value = 2.697141575944443
$ws.Range("C23").Value = 0.1309910452292513
$ws.Range("D23").Value = 0.1117300191871209
$ws.Range("E23").Value = 0.03520330952292605
$ws.Range("G23").Value = 0.002646728475342329
$ws.Range("K23").Value = 2.289549087478747
$ws.Range("L23").Value = 0.2142442327132414
$ws.Range("M23").Value = 0.4955267563879744
$ws.Range("N23").Value = 6.386376519568728
$ws.Range("B24").Value = 2.585295502065662
$ws.Range("C24").Value = 0.114365307667839
$ws.Range("D24").Value = 0.09899564689476392
$ws.Range("E24").Value = 0.03530288926033487
$ws.Range("G24").Value = 0.002655329704040071
$ws.Range("K24").Value = 2.168218424606721
$ws.Range("L24").Value = 0.2091967306545115
$ws.Range("M24").Value = 0.477343206845191
$ws.Range("N24").Value = 6.233050847520417
$ws.Range("B25").Value = 2.47003685507957
$ws.Range("C25").Value = 0.09662033278863191
$ws.Range("D25").Value = 0.0854443229963664
$ws.Range("E25").Value = 0.03543199744222081
$ws.Range("G25").Value = 0.002665281392464486
$ws.Range("K25").Value = 2.042060530952313
$ws.Range("L25").Value = 0.2041402851651668
$ws.Range("M25").Value = 0.4587013185561943
$ws.Range("N25").Value = 6.0694228687573
